$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Removal")
Write-Host $ws.Name
